# Participant Bulk Import Excel Format - edit script
# Applies changes described by the commit:
# "Added PhpSpreadsheet Library and changed the Participant Bulk Import to use PhpSpreadsheet"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Participant LIST")

# 1. Rename the sheet from "Participant LIST" to "Participant List"
$ws.Name = "Participant List"

# 2. Update the Password header note (Q1) text
$ws.Range("Q1").Value = "Password`n- if left blank, the default password is ept1@)(*&^"

# 3. Set K3 to a country value (Afghanistan) - previously blank
$ws.Range("K3").Value = "Afghanistan"

# 4. Adjust row 1 height
$ws.Rows.Item(1).RowHeight = 68

# 5. Update the view: top-left cell and active selection
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K3").Select()
